# OPC UA username and password parameters changed
# Rename the "Username" / "Password" header columns to the fully-qualified
# OPC UA authentication parameter names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "OpcAuthenticationUsername"
$ws.Range("E1").Value = "OpcAuthenticationPassword"

$ws.Range("E2").Select()
